$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4   = 8.853799999999996
    6   = 5.807
    7   = 5.043400000000002
    16  = 6.803799999999999
    20  = 9.413099999999991
    28  = 5.832799999999997
    29  = 5.153200000000004
    32  = 7.093299999999996
    40  = 9.207099999999995
    46  = 5.878099999999996
    51  = 5.6876
    52  = 5.2331
    57  = 5.464899999999997
    59  = 4.890299999999999
    62  = 5.606300000000001
    66  = 5.554199999999996
    73  = 8.797600000000001
    74  = 9.183399999999995
    92  = 4.781899999999996
    100 = 5.301899999999997
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}
